$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the "Python Tracker Connection" activity to the new BronchoVision Viewers entry
#    and record 4 hours against it.
$ws.Range("B42").Value = "* BronchoVision Viewers"
$ws.Range("C42").Value = 4

# 2. Remove the (empty) "Meetings @Parsiss" row - shifts Total Hours / @Parsiss / @Home rows up
$ws.Rows("43:43").Delete()

# 3. Update the @Home hours entry to balance the new total (9 -> 8 hours)
$ws.Range("D45").Value = 8

# 4. Leave the selection where the author ended up editing
[void]$ws.Range("C43").Select()
